$wb = $excel.ActiveWorkbook

# ==================== Sheet: ALC ====================
$ws = $wb.Worksheets.Item("ALC")

# Row 132
$ws.Range("H132").Value = 14072.305
$ws.Range("I132").Value = 5783
$ws.Range("K132").Value = 17349
$ws.Range("M132").Value = -14819

# Row 137
$ws.Range("H137").Value = 9283390
$ws.Range("I137").Value = 557728.3
$ws.Range("J137").Value = 17549806
$ws.Range("K137").Value = 1673184.9
$ws.Range("L137").Value = 52649418
$ws.Range("M137").Value = -1670634.9
$ws.Range("N137").Value = -52654518

# Row 138
$ws.Range("H138").Value = 2005.3529
$ws.Range("I138").Value = 973.75
$ws.Range("J138").Value = 2512.1052
$ws.Range("K138").Value = 2921.25
$ws.Range("L138").Value = 7536.3156
$ws.Range("M138").Value = 2218.75
$ws.Range("N138").Value = -17816.3156

# Row 141
$ws.Range("H141").Value = 4214.9473
$ws.Range("I141").Value = 4214.9473
$ws.Range("K141").Value = 12644.8419
$ws.Range("M141").Value = -7464.841899999999

# ==================== Sheet: ARM ====================
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 18301.48
$ws.Range("I32").Value = 20437.71
$ws.Range("K32").Value = 20437.71
$ws.Range("M32").Value = -20150.71

# Row 63
$ws.Range("H63").Value = 2362.6
$ws.Range("J63").Value = 2377.75
$ws.Range("L63").Value = 2377.75
$ws.Range("N63").Value = -3749.75

# Row 66
$ws.Range("H66").Value = 2362.6
$ws.Range("J66").Value = 2377.75
$ws.Range("L66").Value = 11888.75
$ws.Range("N66").Value = -18752.75

# Row 88
$ws.Range("H88").Value = 251499.5
$ws.Range("I88").Value = 2999
$ws.Range("J88").Value = 500000
$ws.Range("K88").Value = 2999
$ws.Range("L88").Value = 500000
$ws.Range("M88").Value = -2593
$ws.Range("N88").Value = -500812

# Row 91
$ws.Range("H91").Value = 251499.5
$ws.Range("I91").Value = 2999
$ws.Range("J91").Value = 500000
$ws.Range("K91").Value = 2999
$ws.Range("L91").Value = 500000
$ws.Range("M91").Value = -1595
$ws.Range("N91").Value = -502808

# Row 122
$ws.Range("H122").Value = 3654.625
$ws.Range("I122").Value = 2704.5625
$ws.Range("J122").Value = 5554.75
$ws.Range("K122").Value = 8113.6875
$ws.Range("L122").Value = 16664.25
$ws.Range("M122").Value = -5663.6875
$ws.Range("N122").Value = -21564.25

# ==================== Sheet: BSM ====================
$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Range("H20").Value = 14809.125
$ws.Range("J20").Value = 16333.333
$ws.Range("L20").Value = 16333.333
$ws.Range("N20").Value = -16827.333

# Row 86
$ws.Range("H86").Value = 2931.75
$ws.Range("I86").Value = 2563.4
$ws.Range("K86").Value = 2563.4
$ws.Range("M86").Value = -1440.4

# Row 89
$ws.Range("H89").Value = 2931.75
$ws.Range("I89").Value = 2563.4
$ws.Range("K89").Value = 12817
$ws.Range("M89").Value = -7201

# Row 99
$ws.Range("H99").Value = 1603979.4
$ws.Range("I99").Value = 2605266.5
$ws.Range("K99").Value = 2605266.5
$ws.Range("M99").Value = -2603768.5

# Row 105
$ws.Range("H105").Value = 3417
$ws.Range("I105").Value = 3759.8
$ws.Range("K105").Value = 3759.8
$ws.Range("M105").Value = -2012.8

# Row 112
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("N112").ClearContents()

# Row 134
$ws.Range("H134").Value = 1988.6364
$ws.Range("I134").Value = 1287.5
$ws.Range("K134").Value = 3862.5
$ws.Range("M134").Value = -1327.5

# ==================== Sheet: CRP ====================
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 5451.7837
$ws.Range("I31").Value = 1964.2273
$ws.Range("J31").Value = 6927.2886
$ws.Range("K31").Value = 1964.2273
$ws.Range("L31").Value = 6927.2886
$ws.Range("M31").Value = -1669.2273
$ws.Range("N31").Value = -7517.2886

# Row 34
$ws.Range("H34").Value = 5451.7837
$ws.Range("I34").Value = 1964.2273
$ws.Range("J34").Value = 6927.2886
$ws.Range("K34").Value = 1964.2273
$ws.Range("L34").Value = 6927.2886
$ws.Range("M34").Value = -1762.2273
$ws.Range("N34").Value = -7331.2886

# Row 58
$ws.Range("H58").Value = 5105.3335
$ws.Range("I58").Value = 2234.1333
$ws.Range("K58").Value = 2234.1333
$ws.Range("M58").Value = -2031.1333

# Row 122
$ws.Range("H122").Value = 3391.5
$ws.Range("I122").Value = 1976.381
$ws.Range("J122").Value = 6093.091
$ws.Range("K122").Value = 5929.143
$ws.Range("L122").Value = 18279.273
$ws.Range("M122").Value = -3479.143
$ws.Range("N122").Value = -23179.273

# Row 131
$ws.Range("H131").Value = 30000000
$ws.Range("J131").Value = 30000000
$ws.Range("L131").Value = 30000000
$ws.Range("N131").Value = -30010080

# Row 134
$ws.Range("H134").Value = 2258.1428
$ws.Range("I134").Value = 1969.579
$ws.Range("K134").Value = 5908.737
$ws.Range("M134").Value = -3373.737

# Row 136
$ws.Range("H136").Value = 5105.3335
$ws.Range("I136").Value = 2234.1333
$ws.Range("K136").Value = 6702.3999
$ws.Range("M136").Value = -4152.3999

# ==================== Sheet: GSM ====================
$ws = $wb.Worksheets.Item("GSM")

# Row 36
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("N36").ClearContents()

# Row 70
$ws.Range("H70").Value = 6498672.5
$ws.Range("I70").Value = 15154848
$ws.Range("J70").Value = 6541.25
$ws.Range("K70").Value = 15154848
$ws.Range("L70").Value = 6541.25
$ws.Range("M70").Value = -15154578
$ws.Range("N70").Value = -7081.25

# Row 73
$ws.Range("H73").Value = 6498672.5
$ws.Range("I73").Value = 15154848
$ws.Range("J73").Value = 6541.25
$ws.Range("K73").Value = 15154848
$ws.Range("L73").Value = 6541.25
$ws.Range("M73").Value = -15153912
$ws.Range("N73").Value = -8413.25

# Row 80
$ws.Range("H80").Value = 1692024.6
$ws.Range("J80").Value = 48332.832
$ws.Range("L80").Value = 48332.832
$ws.Range("N80").Value = -50328.832

# Row 83
$ws.Range("H83").Value = 1692024.6
$ws.Range("J83").Value = 48332.832
$ws.Range("L83").Value = 241664.16
$ws.Range("N83").Value = -251648.16

# Row 102
$ws.Range("H102").Value = 18526336
$ws.Range("I102").Value = 22734684
$ws.Range("K102").Value = 22734684
$ws.Range("M102").Value = -22733062

# Row 113
$ws.Range("H113").Value = 1997.5
$ws.Range("J113").Value = 1997.5
$ws.Range("L113").Value = 1997.5
$ws.Range("N113").Value = -6337.5

# Row 122
$ws.Range("H122").Value = 463687.88
$ws.Range("I122").Value = 1102938
$ws.Range("K122").Value = 3308814
$ws.Range("M122").Value = -3306364

# Row 123
$ws.Range("H123").Value = 49819.09
$ws.Range("J123").Value = 49819.09
$ws.Range("L123").Value = 49819.09
$ws.Range("N123").Value = -54719.09

# Row 130
$ws.Range("H130").Value = 79998.336
$ws.Range("J130").Value = 79998.336
$ws.Range("L130").Value = 79998.336
$ws.Range("N130").Value = -90038.336

# Row 131
$ws.Range("H131").Value = 100000
$ws.Range("J131").Value = 100000
$ws.Range("L131").Value = 100000
$ws.Range("N131").Value = -110080

# Row 132
$ws.Range("H132").Value = 62818.617
$ws.Range("I132").Value = 86135.53999999999
$ws.Range("J132").Value = 6858
$ws.Range("K132").Value = 258406.62
$ws.Range("L132").Value = 20574
$ws.Range("M132").Value = -255876.62
$ws.Range("N132").Value = -25634

# ==================== Sheet: LTW ====================
$ws = $wb.Worksheets.Item("LTW")

# Row 6
$ws.Range("H6").Value = 79798.5
$ws.Range("J6").Value = 79798.5
$ws.Range("L6").Value = 79798.5
$ws.Range("N6").Value = -80022.5

# Row 38
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("N38").ClearContents()

# Row 40
$ws.Range("H40").Value = 66668450
$ws.Range("I40").Value = 2226
$ws.Range("K40").Value = 2226
$ws.Range("M40").Value = -2090

# Row 100
$ws.Range("H100").Value = 7510.231
$ws.Range("J100").Value = 9241.625
$ws.Range("L100").Value = 9241.625
$ws.Range("N100").Value = -10323.625

# Row 102
$ws.Range("H102").Value = 65000
$ws.Range("J102").Value = 65000
$ws.Range("L102").Value = 65000
$ws.Range("N102").Value = -71490

# Row 117
$ws.Range("H117").Value = 49500
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# Row 122
$ws.Range("H122").Value = 333338000
$ws.Range("J122").Value = 7000
$ws.Range("L122").Value = 21000
$ws.Range("N122").Value = -25900

# Row 123
$ws.Range("H123").Value = 58844.5
$ws.Range("J123").Value = 58844.5
$ws.Range("L123").Value = 58844.5
$ws.Range("N123").Value = -68644.5

# Row 131
$ws.Range("H131").Value = 66532.336
$ws.Range("J131").Value = 66532.336
$ws.Range("L131").Value = 66532.336
$ws.Range("N131").Value = -76612.336

# Row 132
$ws.Range("H132").Value = 4765.269
$ws.Range("I132").Value = 5119.25
$ws.Range("K132").Value = 15357.75
$ws.Range("M132").Value = -12827.75

# ==================== Sheet: WVR ====================
$ws = $wb.Worksheets.Item("WVR")

# Row 15
$ws.Range("H15").Value = 181199.6
$ws.Range("J15").Value = 181199.6
$ws.Range("L15").Value = 181199.6
$ws.Range("N15").Value = -181775.6

# Row 100
$ws.Range("H100").Value = 1080234.8
$ws.Range("I100").Value = 1569432.6
$ws.Range("K100").Value = 3138865.2
$ws.Range("M100").Value = -3138324.2

# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("N102").ClearContents()

# Row 127
$ws.Range("H127").Value = 58800
$ws.Range("J127").Value = 58800
$ws.Range("L127").Value = 58800
$ws.Range("N127").Value = -68720

# Row 129
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("N129").ClearContents()

# Row 132
$ws.Range("H132").Value = 38468520
$ws.Range("I132").Value = 1797.4546
$ws.Range("K132").Value = 5392.3638
$ws.Range("M132").Value = -2862.3638
